# Updated figures for .doc
#
# 1) Rho_Intercept Table: rename the header cell A1 from
#    "Rho_Intercept Name" to "Rho_Intercept Notation".
# 2) Add a new "Apportionment Table" worksheet (as the last tab) describing
#    the Rho_Intercept notation used for spatio-temporal variation in
#    apportionment, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# --- 1) Rename the Rho_Intercept Table header -----------------------------
$wsRho = $wb.Worksheets.Item("Rho_Intercept Table")
$wsRho.Range("A1").Value = "Rho_Intercept Notation"

# --- 2) Add the new Apportionment Table sheet at the end of the workbook --
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsApp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsApp.Name = "Apportionment Table"

# Fill column A (Notation) first
$wsApp.Range("A1").Value = "Notation"
$wsApp.Range("A2").Value = "RW + IaY"
$wsApp.Range("A3").Value = "AR + IaY"
$wsApp.Range("A4").Value = "RW + AR"
$wsApp.Range("A5").Value = "AR + RW"

# Then column C (Spatio-temporal Random Effects)
$wsApp.Range("C1").Value = "Spatio-temporal Random Effects"
$wsApp.Range("C2").Value = "Independent among years"
$wsApp.Range("C3").Value = "Independent among years"
$wsApp.Range("C4").Value = "Autoregressive (lag-1)"
$wsApp.Range("C5").Value = "Random walk"

# Then column B (Intercepts)
$wsApp.Range("B1").Value = "Intercepts"
$wsApp.Range("B2").Value = "Random walk"
$wsApp.Range("B3").Value = "Autoregressive (lag-1)"
$wsApp.Range("B4").Value = "Random walk"
$wsApp.Range("B5").Value = "Autoregressive (lag-1)"

# Bold the header row, matching the other tables in the workbook
$wsApp.Range("A1:C1").Font.Bold = $true

# Widen the notation / intercept columns so the labels aren't truncated
$wsApp.Columns.Item(1).ColumnWidth = 20.497395833333332
$wsApp.Columns.Item(2).ColumnWidth = 18.497395833333332

$wsApp.Range("A1:C5").Select()
